$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

# Column D (Price) updates
Set-TextValue "D2" "274.22"
Set-TextValue "D3" "22.96"
Set-TextValue "D4" "6.344"
Set-TextValue "D5" "0.06242"
Set-TextValue "D6" "3.655"
Set-TextValue "D7" "6.723"
Set-TextValue "D8" "1.368"
Set-TextValue "D9" "0.8314"
Set-TextValue "D10" "0.01376"
Set-TextValue "D11" "0.1634"
Set-TextValue "D12" "0.08307"
Set-TextValue "D13" "0.03361"
Set-TextValue "D14" "0.03102"
Set-TextValue "D15" "0.09312"
Set-TextValue "D16" "3.876"
Set-TextValue "D17" "0.001646"
Set-TextValue "D18" "0.04776"
Set-TextValue "D19" "0.006389"
Set-TextValue "D20" "0.005562"
Set-TextValue "E20" "19HotbitTokenHTBWorstin24h"
Set-TextValue "D21" "0.001088"
Set-TextValue "D23" "3.727"
Set-TextValue "D27" "0.0002679"
Set-TextValue "D40" "0.04700"
Set-TextValue "D41" "0.007025"
Set-TextValue "D42" "0.1164"
Set-TextValue "D43" "0.003599"
Set-TextValue "E43" "42CEJICEJI"
Set-TextValue "D44" "0.01191"
Set-TextValue "D45" "0.00006257"
Set-TextValue "D47" "0.8996"
Set-TextValue "D48" "0.03163"
